# Auto-applies the cell-value updates from the cryptos.xlsx refresh
# (GitHub Actions data refresh commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '29.653.02'
$ws.Range("E2").Value = '  +4.94%  '
$ws.Range("D3").Value = "'" + '1.920.05'
$ws.Range("E3").Value = '  +3.39%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'" + '334.60'
$ws.Range("E5").Value = '  +1.76%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  +2.42%  '
$ws.Range("E8").Value = '  +4.57%  '
$ws.Range("E9").Value = '  +1.24%  '
$ws.Range("D10").Value = "'" + '0.08039'
$ws.Range("E10").Value = '  +3.01%  '
$ws.Range("D11").Value = "'" + '1.014'
$ws.Range("E11").Value = '  +3.28%  '
$ws.Range("D12").Value = "'" + '22.36'
$ws.Range("E12").Value = '  +5.10%  '
$ws.Range("D13").Value = "'" + '1.900.86'
$ws.Range("E13").Value = '  +2.82%  '
$ws.Range("D14").Value = "'" + '5.994'
$ws.Range("E14").Value = '  +3.36%  '
$ws.Range("D15").Value = "'" + '7.175'
$ws.Range("E15").Value = '  +3.14%  '
$ws.Range("D16").Value = "'" + '90.17'
$ws.Range("E16").Value = '  +2.47%  '
$ws.Range("D17").Value = "'" + '1.002'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("E18").Value = '  +1.89%  '
$ws.Range("D19").Value = "'" + '0.06587'
$ws.Range("E19").Value = '  +0.91%  '
$ws.Range("D20").Value = "'" + '17.86'
$ws.Range("E20").Value = '  +4.69%  '
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").Value = "'" + '29.625.20'
$ws.Range("E22").Value = '  +4.82%  '
$ws.Range("D23").Value = "'" + '5.580'
$ws.Range("E23").Value = '  +5.64%  '
$ws.Range("E24").Value = '  +8.03%  '
$ws.Range("D25").Value = "'" + '2.214'
$ws.Range("E25").Value = '  -1.70%  '
$ws.Range("D26").Value = "'" + '2.159.89'
$ws.Range("E26").Value = '  +4.19%  '
$ws.Range("D27").Value = "'" + '156.00'
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("E28").Value = '  +3.46%  '
$ws.Range("D29").Value = "'" + '2.145'
$ws.Range("E29").Value = '  +4.39%  '
$ws.Range("E30").Value = '  +8.66%  '
$ws.Range("D31").Value = "'" + '117.60'
$ws.Range("E31").Value = '  +1.04%  '
$ws.Range("D32").Value = "'" + '1.065'
$ws.Range("E32").Value = '  +13.77%  '
$ws.Range("D33").Value = "'" + '0.09481'
$ws.Range("E33").Value = '  +2.31%  '
$ws.Range("D34").Value = "'" + '1.433'
$ws.Range("E34").Value = '  +4.07%  '
$ws.Range("D35").Value = "'" + '3.579'
$ws.Range("E35").Value = '  -0.57%  '
$ws.Range("E36").Value = '  +4.09%  '
$ws.Range("D37").Value = "'" + '0.06131'
$ws.Range("E37").Value = '  +2.00%  '
$ws.Range("E38").Value = '  +2.74%  '
$ws.Range("D39").Value = "'" + '8.422'
$ws.Range("E39").Value = '  +2.22%  '
$ws.Range("E40").Value = '  +1.67%  '
$ws.Range("D41").Value = "'" + '0.5901'
$ws.Range("E41").Value = '  +3.72%  '
$ws.Range("E42").Value = '  +2.96%  '
$ws.Range("D43").Value = "'" + '10.23'
$ws.Range("E43").Value = '  +2.70%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = "'" + '2.368'
$ws.Range("E44").Value = '  +2.73%  '
$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = "'" + '1.261'
$ws.Range("E45").Value = '  +0.22%  '
$ws.Range("D46").Value = "'" + '0.07515'
$ws.Range("E46").Value = '  +4.99%  '
$ws.Range("D47").Value = "'" + '12.24'
$ws.Range("E47").Value = '  +3.73%  '
$ws.Range("D48").Value = "'" + '0.5568'
$ws.Range("E48").Value = '  +3.30%  '
$ws.Range("D49").Value = "'" + '1.935'
$ws.Range("E49").Value = '  +3.79%  '
$ws.Range("D50").Value = "'" + '112.94'
$ws.Range("E50").Value = '  +2.86%  '
$ws.Range("D51").Value = "'" + '0.3001'
$ws.Range("E51").Value = '  +15.15%  '
